$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 311, pushing existing rows 311-332 down to 312-333
$ws.Rows(311).Insert()

# Populate the newly inserted row 311 with the new weekly data point
$ws.Range("A311").Value = 4
$ws.Range("B311").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C311").Value = "Los Lagos"
$ws.Range("D311").Value = 44931
$ws.Range("D311").NumberFormat = $ws.Range("D312").NumberFormat
$ws.Range("E311").Value = 10
$ws.Range("F311").Value = 100112032
$ws.Range("G311").Value = "Zapallo italiano"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 150
$ws.Range("K311").Value = 12000
$ws.Range("L311").Value = 12000
$ws.Range("M311").Value = 12000
$ws.Range("N311").Value = "`$/caja 50 unidades"
$ws.Range("O311").Value = "Región de O'Higgins"
$ws.Range("P311").Value = 240
$ws.Range("Q311").Value = 50
$ws.Range("R311").Value = "Hortaliza"
